$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ToTitleCase($s) {
    $parts = $s -split "_"
    $result = ""
    for ($j = 0; $j -lt $parts.Count; $j++) {
        $p = $parts[$j]
        if ($j -gt 0) {
            $result = $result + "_"
        }
        $first = $p.Substring(0,1)
        $rest = $p.Substring(1).ToLower()
        $result = $result + $first + $rest
    }
    return $result
}

$lastRow = $ws.UsedRange.Rows.Count()

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $val = $cell.Value()
    $newVal = ToTitleCase($val)
    $cell.Value = $newVal
}
